# Updates the cryptos price/volume snapshot (GitHub Actions data refresh).
# Price cells that look numeric are written with a leading apostrophe so
# Excel stores them as text (matching the sheet's original inlineStr
# formatting, e.g. "581.72" / "0.0250") instead of silently coercing them
# to numbers and dropping significant trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.058.50'
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('D3').Value = '2.646.12'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''581.72'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').Value = '''156.56'
$ws.Range('E6').Value = '  -0.35%  '
$ws.Range('D7').Value = '''0.627'
$ws.Range('E7').Value = '  -2.79%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = '2.644.41'
$ws.Range('E9').Value = '  +0.36%  '
$ws.Range('E10').Value = '  -3.21%  '
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('E12').Value = '  -1.51%  '
$ws.Range('E13').Value = '  +0.95%  '
$ws.Range('D14').Value = '''28.65'
$ws.Range('E14').Value = '  -0.28%  '
$ws.Range('D15').Value = '3.125.52'
$ws.Range('E15').Value = '  +0.39%  '
$ws.Range('D16').Value = '''0.0000185'
$ws.Range('E16').Value = '  -1.03%  '
$ws.Range('D17').Value = '63.970.79'
$ws.Range('D18').Value = '2.645.62'
$ws.Range('E18').Value = '  +0.51%  '
$ws.Range('D19').Value = '''12.23'
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('D20').Value = '''7.76'
$ws.Range('E20').Value = '  +3.75%  '
$ws.Range('D21').Value = '''4.54'
$ws.Range('E21').Value = '  -3.42%  '
$ws.Range('D22').Value = '''346.94'
$ws.Range('E22').Value = '  -0.32%  '
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('D24').Value = '''68.40'
$ws.Range('E24').Value = '  +0.17%  '
$ws.Range('E25').Value = '  +4.40%  '
$ws.Range('E26').Value = '  +0.87%  '
$ws.Range('D27').Value = '''9.38'
$ws.Range('E27').Value = '  -0.84%  '
$ws.Range('B28').Value = 'Fetch.AI'
$ws.Range('C28').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D28').Value = '''1.63'
$ws.Range('E28').Value = '  +2.35%  '
$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D29').Value = '''589.16'
$ws.Range('E29').Value = '  -1.06%  '
$ws.Range('D30').Value = '''8.23'
$ws.Range('E30').Value = '  +2.58%  '
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('D32').Value = '''0.999'
$ws.Range('E32').Value = '  -0.18%  '
$ws.Range('D33').Value = '''2.07'
$ws.Range('E33').Value = '  -0.36%  '
$ws.Range('D34').Value = '''1.75'
$ws.Range('E34').Value = '  +0.36%  '
$ws.Range('D35').Value = '''6.66'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').Value = '''5.50'
$ws.Range('E36').Value = '  +2.89%  '
$ws.Range('D38').Value = '''19.81'
$ws.Range('E38').Value = '  -1.29%  '
$ws.Range('D39').Value = '''1.00'
$ws.Range('E39').Value = '  +0.09%  '
$ws.Range('D40').Value = '''1.92'
$ws.Range('E40').Value = '  -1.05%  '
$ws.Range('D41').Value = '''151.74'
$ws.Range('E41').Value = '  -0.65%  '
$ws.Range('E42').Value = '  +7.41%  '
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Value = '''41.91'
$ws.Range('E44').Value = '  +0.25%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '''164.04'
$ws.Range('E45').Value = '  +3.31%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = '''24.56'
$ws.Range('E46').Value = '  +4.55%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').Value = '''3.93'
$ws.Range('E47').Value = '  -2.38%  '
$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').Value = '''0.0593'
$ws.Range('E48').Value = '  -2.14%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '''0.636'
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').Value = '''0.100'
$ws.Range('E50').Value = '  -1.96%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').Value = '''0.0250'
$ws.Range('E51').Value = '  -2.28%  '
